$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2625
$ws.Range("I12").Value = 250
$ws.Range("J12").Value = 5000
$ws.Range("K12").Value = 250
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = -80
$ws.Range("N12").Value = -5340

$ws.Range("H33").Value = 774.1111
$ws.Range("I33").Value = 776.08
$ws.Range("J33").Value = 749.5
$ws.Range("K33").Value = 776.08
$ws.Range("L33").Value = 749.5
$ws.Range("M33").Value = -547.08
$ws.Range("N33").Value = -1207.5

$ws.Range("H132").Value = 5483.6113
$ws.Range("I132").Value = 7061.6113
$ws.Range("J132").Value = 3905.611
$ws.Range("K132").Value = 21184.8339
$ws.Range("L132").Value = 11716.833
$ws.Range("M132").Value = -18654.8339
$ws.Range("N132").Value = -16776.833

$ws.Range("H138").Value = 1571.5625
$ws.Range("I138").Value = 1009.6667
$ws.Range("J138").Value = 10000
$ws.Range("K138").Value = 3029.0001
$ws.Range("L138").Value = 30000
$ws.Range("M138").Value = 2110.9999
$ws.Range("N138").Value = -40280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11885.676
$ws.Range("I32").Value = 3208.111
$ws.Range("J32").Value = 26904.54
$ws.Range("K32").Value = 3208.111
$ws.Range("L32").Value = 26904.54
$ws.Range("M32").Value = -2921.111
$ws.Range("N32").Value = -27478.54

$ws.Range("H61").Value = 1316
$ws.Range("I61").Value = 1181.5
$ws.Range("J61").Value = 2033.3334
$ws.Range("K61").Value = 1181.5
$ws.Range("L61").Value = 2033.3334
$ws.Range("M61").Value = -969.5
$ws.Range("N61").Value = -2457.3334

$ws.Range("H64").Value = 16660
$ws.Range("J64").Value = 16660
$ws.Range("L64").Value = 16660
$ws.Range("N64").Value = -17156

$ws.Range("H67").Value = 16660
$ws.Range("J67").Value = 16660
$ws.Range("L67").Value = 16660
$ws.Range("N67").Value = -18376

$ws.Range("H74").Value = 9617252
$ws.Range("I74").Value = 15626039
$ws.Range("J74").Value = 3194
$ws.Range("K74").Value = 15626039
$ws.Range("L74").Value = 3194
$ws.Range("M74").Value = -15625165
$ws.Range("N74").Value = -4942

$ws.Range("H77").Value = 9617252
$ws.Range("I77").Value = 15626039
$ws.Range("J77").Value = 3194
$ws.Range("K77").Value = 78130195
$ws.Range("L77").Value = 15970
$ws.Range("M77").Value = -78125827
$ws.Range("N77").Value = -24706

$ws.Range("H132").Value = 1356.1923
$ws.Range("I132").Value = 1231.4791
$ws.Range("J132").Value = 2852.75
$ws.Range("K132").Value = 3694.4373
$ws.Range("L132").Value = 8558.25
$ws.Range("M132").Value = -1164.4373
$ws.Range("N132").Value = -13618.25

$ws.Range("H136").Value = 1316
$ws.Range("I136").Value = 1181.5
$ws.Range("J136").Value = 2033.3334
$ws.Range("K136").Value = 3544.5
$ws.Range("L136").Value = 6100.0002
$ws.Range("M136").Value = -994.5
$ws.Range("N136").Value = -11200.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 36333.332
$ws.Range("J62").Value = 36333.332
$ws.Range("L62").Value = 36333.332
$ws.Range("N62").Value = -37705.332

$ws.Range("H65").Value = 36333.332
$ws.Range("J65").Value = 36333.332
$ws.Range("L65").Value = 108999.996
$ws.Range("N65").Value = -115863.996

$ws.Range("H86").Value = 10528050
$ws.Range("I86").Value = 12501497
$ws.Range("K86").Value = 12501497
$ws.Range("M86").Value = -12500374

$ws.Range("H89").Value = 10528050
$ws.Range("I89").Value = 12501497
$ws.Range("K89").Value = 62507485
$ws.Range("M89").Value = -62501869

$ws.Range("H134").Value = 1997.35
$ws.Range("I134").Value = 1556.4667
$ws.Range("J134").Value = 3320
$ws.Range("K134").Value = 4669.4001
$ws.Range("L134").Value = 9960
$ws.Range("M134").Value = -2134.4001
$ws.Range("N134").Value = -15030

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1046.9
$ws.Range("I7").Value = 29.8
$ws.Range("J7").Value = 2064
$ws.Range("K7").Value = 29.8
$ws.Range("L7").Value = 2064
$ws.Range("M7").Value = 83.2
$ws.Range("N7").Value = -2290

$ws.Range("H31").Value = 4188854.5
$ws.Range("I31").Value = 6700963.5
$ws.Range("J31").Value = 2005.5555
$ws.Range("K31").Value = 6700963.5
$ws.Range("L31").Value = 2005.5555
$ws.Range("M31").Value = -6700668.5
$ws.Range("N31").Value = -2595.5555

$ws.Range("H34").Value = 4188854.5
$ws.Range("I34").Value = 6700963.5
$ws.Range("J34").Value = 2005.5555
$ws.Range("K34").Value = 6700963.5
$ws.Range("L34").Value = 2005.5555
$ws.Range("M34").Value = -6700761.5
$ws.Range("N34").Value = -2409.5555

$ws.Range("H58").Value = 1730.9
$ws.Range("I58").Value = 794.5
$ws.Range("J58").Value = 2355.1667
$ws.Range("K58").Value = 794.5
$ws.Range("L58").Value = 2355.1667
$ws.Range("M58").Value = -591.5
$ws.Range("N58").Value = -2761.1667

$ws.Range("H99").Value = 3577524.8
$ws.Range("I99").Value = 4470285
$ws.Range("J99").Value = 6485
$ws.Range("K99").Value = 4470285
$ws.Range("L99").Value = 6485
$ws.Range("M99").Value = -4468787
$ws.Range("N99").Value = -9481

$ws.Range("H126").Value = 3577524.8
$ws.Range("I126").Value = 4470285
$ws.Range("J126").Value = 6485
$ws.Range("K126").Value = 13410855
$ws.Range("L126").Value = 19455
$ws.Range("M126").Value = -13408385
$ws.Range("N126").Value = -24395

$ws.Range("H132").Value = 1261.4166
$ws.Range("I132").Value = 1012.5
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 3037.5
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -507.5
$ws.Range("N132").Value = -17058.5

$ws.Range("H134").Value = 4354.8335
$ws.Range("I134").Value = 9724.5
$ws.Range("J134").Value = 1670
$ws.Range("K134").Value = 29173.5
$ws.Range("L134").Value = 5010
$ws.Range("M134").Value = -26638.5
$ws.Range("N134").Value = -10080

$ws.Range("H136").Value = 1730.9
$ws.Range("I136").Value = 794.5
$ws.Range("J136").Value = 2355.1667
$ws.Range("K136").Value = 2383.5
$ws.Range("L136").Value = 7065.500100000001
$ws.Range("M136").Value = 166.5
$ws.Range("N136").Value = -12165.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4158.407
$ws.Range("I122").Value = 3313.5
$ws.Range("J122").Value = 6572.4287
$ws.Range("K122").Value = 9940.5
$ws.Range("L122").Value = 19717.2861
$ws.Range("M122").Value = -7490.5
$ws.Range("N122").Value = -24617.2861

$ws.Range("H132").Value = 2499.48
$ws.Range("I132").Value = 1793.5
$ws.Range("J132").Value = 3398
$ws.Range("K132").Value = 5380.5
$ws.Range("L132").Value = 10194
$ws.Range("M132").Value = -2850.5
$ws.Range("N132").Value = -15254

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 15251130
$ws.Range("I132").Value = 34735660
$ws.Range("J132").Value = 2367.5217
$ws.Range("K132").Value = 104206980
$ws.Range("L132").Value = 7102.5651
$ws.Range("M132").Value = -104204450
$ws.Range("N132").Value = -12162.5651

$ws.Range("H136").Value = 5655.7334
$ws.Range("I136").Value = 10685.538
$ws.Range("J136").Value = 1809.4117
$ws.Range("K136").Value = 32056.614
$ws.Range("L136").Value = 5428.2351
$ws.Range("M136").Value = -29506.614
$ws.Range("N136").Value = -10528.2351

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4654
$ws.Range("I132").Value = 5399.5
$ws.Range("J132").Value = 4355.8
$ws.Range("K132").Value = 16198.5
$ws.Range("L132").Value = 13067.4
$ws.Range("M132").Value = -13668.5
$ws.Range("N132").Value = -18127.4

$ws.Range("H136").Value = 1565.5588
$ws.Range("I136").Value = 1041.4706
$ws.Range("J136").Value = 2089.647
$ws.Range("K136").Value = 3124.4118
$ws.Range("L136").Value = 6268.941
$ws.Range("M136").Value = -574.4118000000003
$ws.Range("N136").Value = -11368.941
